$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - Control 46
$ws.Range("D2").Value = 0.003489333874770055
$ws.Range("E2").Value = 0.003489333874770055

# Row 3 - Control 28
$ws.Range("D3").Value = 0.9962159613032793
$ws.Range("E3").Value = 0.9962159613032793

# Row 4 - Control 13
$ws.Range("D4").Value = 0.9995512664259458
$ws.Range("E4").Value = 0.9995512664259458

# Row 5 - Control 50
$ws.Range("D5").Value = 0.0000000000000000000000000005357241081620243
$ws.Range("E5").Value = 0.0000000000000000000000000005357241081620243

# Row 6 - Control 51
$ws.Range("D6").Value = 0.00000002578622050009275
$ws.Range("E6").Value = 0.00000002578622050009275

# Row 7 - MDD 41
$ws.Range("D7").Value = 0.00002660187187741686
$ws.Range("E7").Value = 0.9999733981281226

# Row 8 - MDD 8
$ws.Range("D8").Value = 0.9999822589184943
$ws.Range("E8").Value = 0.00001774108150565112

# Row 9 - MDD 15
$ws.Range("C9").Value = $true
$ws.Range("D9").Value = 0.9195485162839092
$ws.Range("E9").Value = 0.08045148371609079

# Row 10 - MDD 16
$ws.Range("D10").Value = 0.9999999992022208
$ws.Range("E10").Value = 0.0000000007977791760538366

# Row 11 - MDD 33
$ws.Range("D11").Value = 0.9999999999992655
$ws.Range("E11").Value = 0.0000000000007345235530920036
$ws.Range("F11").Value = 2.390795946121216
$ws.Range("G11").Value = 0.7
